$wb = $excel.ActiveWorkbook

# --- Sheet 1: Weekly Quantity ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")

# Update B17 value
$ws1.Cells.Item(17, 2).Value = 758

# Delete old rows 30-33 (4 rows), shifting rows 34-54 up to become rows 30-50
$ws1.Range("A30:B33").EntireRow.Delete()

# --- Sheet 2: Monthly Trend ---
$ws2 = $wb.Worksheets.Item("Monthly Trend")

# Update B7 value
$ws2.Cells.Item(7, 2).Value = 1010

# Delete old row 13, shifting rows 14-20 up to become rows 13-19
$ws2.Range("A13:B13").EntireRow.Delete()
